# Update countries & provincias Spain
# Refresh the "Pais" sheet: new case totals caused several countries to
# change rank (they are listed ordered by "Casos totales" descending), so
# both the country label and the statistics of the affected rows are
# rewritten here; the timestamp caption in A1 is updated too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'Datos actualizados a 3 de Abril de 2020 a las 01:20'

$ws.Range("A4").Value = 'Estados Unidos'
$ws.Range("B4").Value = 243970
$ws.Range("C4").Value = 28967
$ws.Range("D4").Value = 10400
$ws.Range("E4").Value = 227687
$ws.Range("F4").Value = 5421
$ws.Range("G4").Value = 781
$ws.Range("H4").Value = 5883

$ws.Range("A17").Value = 'Austria'
$ws.Range("B17").Value = 11129
$ws.Range("C17").Value = 418
$ws.Range("D17").Value = 1749
$ws.Range("E17").Value = 9222
$ws.Range("F17").Value = 227
$ws.Range("G17").Value = 12
$ws.Range("H17").Value = 158

$ws.Range("A20").Value = 'Brasil'
$ws.Range("B20").Value = 8044
$ws.Range("C20").Value = 1164
$ws.Range("D20").Value = 127
$ws.Range("E20").Value = 7593
$ws.Range("F20").Value = 296
$ws.Range("G20").Value = 82
$ws.Range("H20").Value = 324

$ws.Range("A24").Value = 'Noruega'
$ws.Range("B24").Value = 5147
$ws.Range("C24").Value = 270
$ws.Range("D24").Value = 32
$ws.Range("E24").Value = 5065
$ws.Range("F24").Value = 96
$ws.Range("G24").Value = 6
$ws.Range("H24").Value = 50

$ws.Range("A25").Value = 'Chequia'
$ws.Range("B25").Value = 3858
$ws.Range("C25").Value = 269
$ws.Range("D25").Value = 67
$ws.Range("E25").Value = 3747
$ws.Range("F25").Value = 72
$ws.Range("G25").Value = 5
$ws.Range("H25").Value = 44

$ws.Range("A26").Value = 'Irlanda'
$ws.Range("B26").Value = 3849
$ws.Range("C26").Value = 402
$ws.Range("D26").Value = 5
$ws.Range("E26").Value = 3746
$ws.Range("F26").Value = 109
$ws.Range("G26").Value = 13
$ws.Range("H26").Value = 98

$ws.Range("A36").Value = 'Japon'
$ws.Range("B36").Value = 2495
$ws.Range("C36").Value = 111
$ws.Range("D36").Value = 472
$ws.Range("E36").Value = 1961
$ws.Range("F36").Value = 60
$ws.Range("G36").Value = 5
$ws.Range("H36").Value = 62

$ws.Range("A37").Value = 'Luxemburgo'
$ws.Range("B37").Value = 2487
$ws.Range("C37").Value = 168
$ws.Range("D37").Value = 80
$ws.Range("E37").Value = 2377
$ws.Range("F37").Value = 31
$ws.Range("G37").Value = 1
$ws.Range("H37").Value = 30

$ws.Range("A38").Value = 'Pakistan'
$ws.Range("B38").Value = 2421
$ws.Range("C38").Value = 303
$ws.Range("D38").Value = 125
$ws.Range("E38").Value = 2262
$ws.Range("F38").Value = 9
$ws.Range("G38").Value = 7
$ws.Range("H38").Value = 34

$ws.Range("A52").Value = 'Argentina'
$ws.Range("B52").Value = 1133
$ws.Range("C52").Value = 0
$ws.Range("D52").Value = 256
$ws.Range("E52").Value = 841
$ws.Range("F52").Value = 0
$ws.Range("G52").Value = 4
$ws.Range("H52").Value = 36

$ws.Range("A65").Value = 'Crucero'
$ws.Range("B65").Value = 712
$ws.Range("C65").Value = 0
$ws.Range("D65").Value = 619
$ws.Range("E65").Value = 82
$ws.Range("F65").Value = 10
$ws.Range("G65").Value = 0
$ws.Range("H65").Value = 11

$ws.Range("A96").Value = 'Vietnam'
$ws.Range("B96").Value = 233
$ws.Range("C96").Value = 15
$ws.Range("D96").Value = 75
$ws.Range("E96").Value = 158
$ws.Range("F96").Value = 3
$ws.Range("G96").Value = 0
$ws.Range("H96").Value = 0

$ws.Range("A97").Value = 'Oman'
$ws.Range("B97").Value = 231
$ws.Range("C97").Value = 21
$ws.Range("D97").Value = 57
$ws.Range("E97").Value = 173
$ws.Range("F97").Value = 3
$ws.Range("G97").Value = 0
$ws.Range("H97").Value = 1

$ws.Range("A112").Value = 'Consejo Danes para los Refugiados'
$ws.Range("B112").Value = 134
$ws.Range("C112").Value = 25
$ws.Range("D112").Value = 3
$ws.Range("E112").Value = 118
$ws.Range("F112").Value = 0
$ws.Range("G112").Value = 4
$ws.Range("H112").Value = 13

$ws.Range("A113").Value = 'Georgia'
$ws.Range("B113").Value = 134
$ws.Range("C113").Value = 17
$ws.Range("D113").Value = 26
$ws.Range("E113").Value = 108
$ws.Range("F113").Value = 6
$ws.Range("G113").Value = 0
$ws.Range("H113").Value = 0

$ws.Range("A114").Value = 'Brunei'
$ws.Range("B114").Value = 133
$ws.Range("C114").Value = 2
$ws.Range("D114").Value = 56
$ws.Range("E114").Value = 76
$ws.Range("F114").Value = 3
$ws.Range("G114").Value = 0
$ws.Range("H114").Value = 1

$ws.Range("A115").Value = 'Guadalupe'
$ws.Range("B115").Value = 128
$ws.Range("C115").Value = 3
$ws.Range("D115").Value = 24
$ws.Range("E115").Value = 98
$ws.Range("F115").Value = 14
$ws.Range("G115").Value = 0
$ws.Range("H115").Value = 6

$ws.Range("A116").Value = 'Bolivia'
$ws.Range("B116").Value = 123
$ws.Range("C116").Value = 8
$ws.Range("D116").Value = 1
$ws.Range("E116").Value = 114
$ws.Range("F116").Value = 3
$ws.Range("G116").Value = 1
$ws.Range("H116").Value = 8

$ws.Range("A149").Value = 'Islas Caimanes'
$ws.Range("B149").Value = 28
$ws.Range("C149").Value = 6
$ws.Range("D149").Value = 0
$ws.Range("E149").Value = 27
$ws.Range("F149").Value = 0
$ws.Range("G149").Value = 0
$ws.Range("H149").Value = 1

$ws.Range("A150").Value = 'Bahamas'
$ws.Range("B150").Value = 24
$ws.Range("C150").Value = 3
$ws.Range("D150").Value = 1
$ws.Range("E150").Value = 22
$ws.Range("F150").Value = 0
$ws.Range("G150").Value = 0
$ws.Range("H150").Value = 1

$ws.Range("A151").Value = 'Eritrea'
$ws.Range("B151").Value = 22
$ws.Range("C151").Value = 7
$ws.Range("D151").Value = 0
$ws.Range("E151").Value = 22
$ws.Range("F151").Value = 0
$ws.Range("G151").Value = 0
$ws.Range("H151").Value = 0

$ws.Range("A152").Value = 'San Martin (Parte Francesa)'
$ws.Range("B152").Value = 22
$ws.Range("C152").Value = 0
$ws.Range("D152").Value = 2
$ws.Range("E152").Value = 19
$ws.Range("F152").Value = 0
$ws.Range("G152").Value = 0
$ws.Range("H152").Value = 1

$ws.Range("A153").Value = 'Congo'
$ws.Range("B153").Value = 22
$ws.Range("C153").Value = 0
$ws.Range("D153").Value = 2
$ws.Range("E153").Value = 18
$ws.Range("F153").Value = 0
$ws.Range("G153").Value = 0
$ws.Range("H153").Value = 2

$ws.Range("A154").Value = 'Gabon'
$ws.Range("B154").Value = 21
$ws.Range("C154").Value = 3
$ws.Range("D154").Value = 0
$ws.Range("E154").Value = 20
$ws.Range("F154").Value = 0
$ws.Range("G154").Value = 0
$ws.Range("H154").Value = 1

$ws.Range("A167").Value = 'Santa Lucia'
$ws.Range("B167").Value = 13
$ws.Range("C167").Value = 0
$ws.Range("D167").Value = 1
$ws.Range("E167").Value = 12
$ws.Range("F167").Value = 0
$ws.Range("G167").Value = 0
$ws.Range("H167").Value = 0

$ws.Range("A168").Value = 'Benin'
$ws.Range("B168").Value = 13
$ws.Range("C168").Value = 0
$ws.Range("D168").Value = 1
$ws.Range("E168").Value = 12
$ws.Range("F168").Value = 0
$ws.Range("G168").Value = 0
$ws.Range("H168").Value = 0

$ws.Range("A172").Value = 'Seychelles'
$ws.Range("B172").Value = 10
$ws.Range("C172").Value = 0
$ws.Range("D172").Value = 0
$ws.Range("E172").Value = 10
$ws.Range("F172").Value = 0
$ws.Range("G172").Value = 0
$ws.Range("H172").Value = 0

$ws.Range("A173").Value = 'Laos'
$ws.Range("B173").Value = 10
$ws.Range("C173").Value = 0
$ws.Range("D173").Value = 0
$ws.Range("E173").Value = 10
$ws.Range("F173").Value = 0
$ws.Range("G173").Value = 0
$ws.Range("H173").Value = 0

$ws.Range("A174").Value = 'Surinam'
$ws.Range("B174").Value = 10
$ws.Range("C174").Value = 0
$ws.Range("D174").Value = 0
$ws.Range("E174").Value = 10
$ws.Range("F174").Value = 0
$ws.Range("G174").Value = 0
$ws.Range("H174").Value = 0

$ws.Range("A175").Value = 'Granada'
$ws.Range("B175").Value = 10
$ws.Range("C175").Value = 1
$ws.Range("D175").Value = 0
$ws.Range("E175").Value = 10
$ws.Range("F175").Value = 0
$ws.Range("G175").Value = 0
$ws.Range("H175").Value = 0

$ws.Range("A176").Value = 'Mozambique'
$ws.Range("B176").Value = 10
$ws.Range("C176").Value = 0
$ws.Range("D176").Value = 0
$ws.Range("E176").Value = 10
$ws.Range("F176").Value = 0
$ws.Range("G176").Value = 0
$ws.Range("H176").Value = 0

$ws.Range("A179").Value = 'San Cristobal y Nieves'
$ws.Range("B179").Value = 9
$ws.Range("C179").Value = 1
$ws.Range("D179").Value = 0
$ws.Range("E179").Value = 9
$ws.Range("F179").Value = 0
$ws.Range("G179").Value = 0
$ws.Range("H179").Value = 0

$ws.Range("A181").Value = 'Guinea-Bisau'
$ws.Range("B181").Value = 9
$ws.Range("C181").Value = 0
$ws.Range("D181").Value = 0
$ws.Range("E181").Value = 9
$ws.Range("F181").Value = 0
$ws.Range("G181").Value = 0
$ws.Range("H181").Value = 0

$ws.Range("A195").Value = 'Nicaragua'
$ws.Range("B195").Value = 5
$ws.Range("C195").Value = 0
$ws.Range("D195").Value = 0
$ws.Range("E195").Value = 4
$ws.Range("F195").Value = 0
$ws.Range("G195").Value = 0
$ws.Range("H195").Value = 1

$ws.Range("A196").Value = 'Somalia'
$ws.Range("B196").Value = 5
$ws.Range("C196").Value = 0
$ws.Range("D196").Value = 1
$ws.Range("E196").Value = 4
$ws.Range("F196").Value = 0
$ws.Range("G196").Value = 0
$ws.Range("H196").Value = 0

$ws.Range("A201").Value = 'Republica de Africa Central'
$ws.Range("B201").Value = 3
$ws.Range("C201").Value = 0
$ws.Range("D201").Value = 0
$ws.Range("E201").Value = 3
$ws.Range("F201").Value = 0
$ws.Range("G201").Value = 0
$ws.Range("H201").Value = 0

$ws.Range("A202").Value = 'Belice'
$ws.Range("B202").Value = 3
$ws.Range("C202").Value = 0
$ws.Range("D202").Value = 0
$ws.Range("E202").Value = 3
$ws.Range("F202").Value = 0
$ws.Range("G202").Value = 0
$ws.Range("H202").Value = 0

$ws.Range("A203").Value = 'Islas Virgenes Britanicas'
$ws.Range("B203").Value = 3
$ws.Range("C203").Value = 0
$ws.Range("D203").Value = 0
$ws.Range("E203").Value = 3
$ws.Range("F203").Value = 0
$ws.Range("G203").Value = 0
$ws.Range("H203").Value = 0

$ws.Range("A209").Value = 'Papua Nueva Guinea'
$ws.Range("B209").Value = 1
$ws.Range("C209").Value = 0
$ws.Range("D209").Value = 0
$ws.Range("E209").Value = 1
$ws.Range("F209").Value = 0
$ws.Range("G209").Value = 0
$ws.Range("H209").Value = 0

$ws.Range("A210").Value = 'Timor Oriental'
$ws.Range("B210").Value = 1
$ws.Range("C210").Value = 0
$ws.Range("D210").Value = 0
$ws.Range("E210").Value = 1
$ws.Range("F210").Value = 0
$ws.Range("G210").Value = 0
$ws.Range("H210").Value = 0
